$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture existing values before they get overwritten ---
$oldD1 = $ws.Range("D1").Value2
$oldD2 = $ws.Range("D2").Value2
$oldD3 = $ws.Range("D3").Value2
$oldD4 = $ws.Range("D4").Value2
$oldD5 = $ws.Range("D5").Value2
$oldD6 = $ws.Range("D6").Value2
$oldD7 = $ws.Range("D7").Value2
$oldD8 = $ws.Range("D8").Value2

$oldC6 = $ws.Range("C6").Value2
$oldC7 = $ws.Range("C7").Value2

# --- Split row 6 prerequisites into "prerequisite" + "recommended" parts ---
$marker6 = "Recommended: "
$idx6 = $oldC6.IndexOf($marker6)
$newC6 = $oldC6.Substring(0, $idx6).TrimEnd()
$newF6 = $oldC6.Substring($idx6 + $marker6.Length)

# --- Remove the word "courses:" from row 7 prerequisites text ---
$newC7 = $oldC7.Replace("following courses: ", "following ")

# --- Move "Terms Typically Offered" column from D to its new home, G ---
$ws.Range("G1").Value = $oldD1
$ws.Range("G2").Value = $oldD2
$ws.Range("G3").Value = $oldD3
$ws.Range("G4").Value = $oldD4
$ws.Range("G5").Value = $oldD5
$ws.Range("G6").Value = $oldD6 + " "
$ws.Range("G7").Value = $oldD7
$ws.Range("G8").Value = $oldD8

# --- New headers for the inserted columns ---
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# --- New column D values (Corequisites) ---
$ws.Range("D2").Value = "NA"
$ws.Range("D3").Value = "NA"
$ws.Range("D4").Value = "NA"
$ws.Range("D5").Value = "NA"
$ws.Range("D6").Value = "NA"
$ws.Range("D7").Value = "NA"
$ws.Range("D8").Value = "NA"

# --- New column E values (Concurrent) ---
$ws.Range("E2").Value = "NA"
$ws.Range("E3").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("E5").Value = "NA"
$ws.Range("E6").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("E8").Value = "NA"

# --- New column F values (Recommended) ---
$ws.Range("F2").Value = "NA"
$ws.Range("F3").Value = "NA"
$ws.Range("F4").Value = "NA"
$ws.Range("F5").Value = "NA"
$ws.Range("F6").Value = $newF6
$ws.Range("F7").Value = "NA"
$ws.Range("F8").Value = "NA"

# --- Updated Prerequisites text (column C) for rows 6 and 7 ---
$ws.Range("C6").Value = $newC6
$ws.Range("C7").Value = $newC7
